$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9808467030525208
$ws.Range("B1").Value = 1.195519924163818
$ws.Range("C1").Value = 1.022421598434448
$ws.Range("D1").Value = 0.9515233039855957
$ws.Range("E1").Value = 0.9894855618476868
